$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.756.90"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.938.39"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "567.41"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").Value = "158.00"
$ws.Range("E6").Value = "  +2.50%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").Value = "2.933.00"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -5.28%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "33.96"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "64.943.19"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "3.421.96"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "2.940.83"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "444.95"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").Value = "13.81"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "7.23"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "82.57"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "12.01"
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("E27").Value = "  -6.78%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "7.95"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").Value = "2.36"
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "0.0₃0993"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "27.13"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "0.976"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").Value = "5.67"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").Value = "49.08"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").Value = "  -5.92%  "
$ws.Range("D40").Value = "43.51"
$ws.Range("E40").Value = "  -3.86%  "
$ws.Range("D41").Value = "0.297"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "382.35"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").Value = "2.723.45"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "132.01"
$ws.Range("E48").Value = "  -1.76%  "

# Row 42/43 swap: Kaspa <-> dogwifhat with updated values
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.82"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.118"
$ws.Range("E43").Value = "  -1.91%  "

# Row 51 replacement: ThetaToken -> InjectiveProtocol
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "23.19"
$ws.Range("E51").Value = "  -0.11%  "
